# Update "last_modified:" / "Date_Modified:" timestamp columns to reflect
# the files being moved to a new repo (ETL portion with postgres moved out).
# These values are raw Excel date serial numbers (days since 1899-12-30).

$wb = $excel.ActiveWorkbook

# Sheet "data-input": column W holds last_modified: for rows 11-13
$wsData = $wb.Worksheets.Item("data-input")
$wsData.Range("W11").Value2 = 43970.45450629989
$wsData.Range("W12").Value2 = 43970.45450760752
$wsData.Range("W13").Value2 = 43970.45450830195

# Sheet "dir_data-input": column E holds Date_Modified: for rows 4, 10, 13
$wsDir = $wb.Worksheets.Item("dir_data-input")
$wsDir.Range("E4").Value2 = 43970.4545054895
$wsDir.Range("E10").Value2 = 43970.45450720243
$wsDir.Range("E13").Value2 = 43970.45450803581
